# "Pagina Menu Inicio Commit"
# The title textbox on slide 1 ("CuadroTexto 5") had its first line
# "Portafolio - Desarrollador Web" retyped. The hyphen became an en-dash
# and a stray tab was left after "D", which also split the run that used
# to hold the whole line into three separate runs (PowerPoint re-split the
# run because the middle portion - "esarrollador" - got flagged by the
# spell checker once it was no longer part of the word "Desarrollador").

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$shape = $slide.Shapes.Item(1)
$titleRange = $shape.TextFrame.TextRange
$firstLine = $titleRange.Paragraphs(1)

# Replace "Portafolio - D" (14 chars) with "Portafolio" + en-dash + " D" + TAB
$enDash = [char]0x2013
$tab = [char]0x09
$newStart = "Portafolio " + $enDash + " D" + $tab
$firstLine.Characters(1, 14).Text = $newStart

# Re-assign the remaining text in place so it becomes its own runs,
# matching the run split seen in the saved file.
$firstLine.Characters(16, 12).Text = "esarrollador"
$firstLine.Characters(28, 4).Text = " Web"
